{"js": "// Docx writer: Use different style for block quotes in notes.\n//\n// Adds a new paragraph style \"Footnote Block Text\" (styleId\n// \"FootnoteBlockText\"), based on \"Footnote Text\" (and followed by\n// \"Footnote Text\"), with the same block-quote spacing/indent formatting\n// as the existing \"Block Text\" style. This gives footnote block quotes\n// their own style so they can later be given a different font size than\n// the rest of the footnote text.\n\n// 1. Create the new paragraph style.\ncontext.document.addStyle(\"Footnote Block Text\", Word.StyleType.paragraph);\nawait context.sync();\n\n// 2. Re-fetch a live handle to the style we just created. (The object\n// returned directly by addStyle() cannot be used to set further\n// properties reliably, so look it up by name after syncing.)\nconst style = context.document.getStyles().getByNameOrNullObject(\"Footnote Block Text\");\nawait context.sync();\n\nif (style.isNullObject) {\n  throw new Error(\"Failed to create the 'Footnote Block Text' style.\");\n}\n\n// 3. basedOn / next paragraph style -> \"Footnote Text\".\nstyle.baseStyle = \"Footnote Text\";\nstyle.nextParagraphStyle = \"Footnote Text\";\n\n// 4. Misc style metadata matching the other custom block-level styles\n// in this stylesheet (e.g. \"Block Text\"): uiPriority=9, unhideWhenUsed,\n// qFormat.\nstyle.priority = 9;\nstyle.unhideWhenUsed = true;\nstyle.quickStyle = true;\n\n// 5. Paragraph formatting: same spacing/indent as \"Block Text\" -\n// spacing before/after = 5pt (100 twips), first line indent = 0,\n// left/right indent = 24pt (480 twips).\nconst pf = style.paragraphFormat;\npf.spaceBefore = 5;\npf.spaceAfter = 5;\npf.firstLineIndent = 0;\npf.leftIndent = 24;\npf.rightIndent = 24;\n\nawait context.sync();\n", "ps1": "# Docx writer: Use different style for block quotes in notes.\n#\n# Adds a new paragraph style \"Footnote Block Text\" (styleId\n# \"FootnoteBlockText\"), based on \"Footnote Text\" (and followed by\n# \"Footnote Text\"), with the same block-quote spacing/indent formatting\n# as the existing \"Block Text\" style. This gives footnote block quotes\n# their own style so they can later be given a different font size than\n# the rest of the footnote text.\n\n$d = $word.ActiveDocument\n\n# 1 = wdStyleTypeParagraph\n$style = $d.Styles.Add(\"Footnote Block Text\", 1)\n\n# basedOn / next paragraph style -> \"Footnote Text\".\n$style.BaseStyle = \"Footnote Text\"\n$style.NextParagraphStyle = \"Footnote Text\"\n\n# Misc style metadata matching the other custom block-level styles in\n# this stylesheet (e.g. \"Block Text\"): uiPriority=9, unhideWhenUsed,\n# qFormat.\n$style.Priority = 9\n$style.UnhideWhenUsed = $true\n$style.QuickStyle = $true\n\n# Paragraph formatting: same spacing/indent as \"Block Text\" - spacing\n# before/after = 5pt (100 twips), first line indent = 0, left/right\n# indent = 24pt (480 twips).\n$pf = $style.ParagraphFormat\n$pf.SpaceBefore = 5\n$pf.SpaceAfter = 5\n$pf.FirstLineIndent = 0\n$pf.LeftIndent = 24\n$pf.RightIndent = 24\n"}
